# The workbook originally had a "datetime" column (column A) that was
# removed from the data set. In the UI this corresponds to: click the
# column A header to select the whole column, then delete it (Home >
# Delete > Delete Sheet Columns, or right-click > Delete). Excel shifts
# every column to the right of A one place to the left (B->A, C->B, ...),
# drops the shared string "datetime" (no longer referenced), and the
# number-format/style that only the datetime cells used becomes unused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire first column (mirrors clicking the "A" column header)
$ws.Columns.Item(1).Select()

# Delete it, shifting everything else left
$ws.Columns.Item(1).Delete()
